$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns are treated as text so values like "1.00" or "2.40"
# are not silently normalized into numbers (matches original inlineStr cells).
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.015.28'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.55%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.301.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.98%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.26'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.90%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.02'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.28%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.63%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.511'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.14%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.41'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.43%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0788'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.71%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.117'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.74%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.84'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.00%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.81'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.75%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.657.92'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.05%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.295.92'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.24%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.785'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.24%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.963.61'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.51%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.67'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.53%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.07%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.12'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.93%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.63'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.57%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '242.28'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.34%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.14'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.63%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.03%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.43'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.35%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.02'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.32%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.11'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.16%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.33'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.46%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.74%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.08'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.16%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.21'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.46%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.01%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.04'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.97%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.75'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.11%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.78'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.58%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.40'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.39%  '

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.10%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.41%  '

# Row 40
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.77'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.85%  '

# Row 41
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.77'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.25%  '

# Row 42
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.111'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.39%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.004.71'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.09%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0285'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.78%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.18'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.41%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.25'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.97%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.36'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.93%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.81'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.89%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.86'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.03%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.524.01'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.05%  '

# Row 51
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.52'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.86%  '
